$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SECOND YEAR")

# Updated grade-distribution counts (number of students per grade bucket).
$ws.Range("E9").Value = 0
$ws.Range("E11").Value = 8
$ws.Range("E12").Value = 2
$ws.Range("E14").Value = 15
$ws.Range("E16").Value = 14
$ws.Range("E18").Value = 13

# Leave the sheet active with the cell that was last edited selected.
$ws.Activate()
$ws.Range("E13").Select()
